$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Fix typo in B31: "צקרה" -> "תקרה"
$ws.Range("B31").Value = "הצעת מקום תקרה אם אין מקום פנוי בזמן השלקוח רוצה או שאין שולחן פנוי בגודל שהוא רוצה "

# 2) Fix typo in E19: "יצבצע" -> "יתבצע"
$ws.Range("E19").Value = "פינוי שולחן דורש שהתשלום יתבצע לפניו "

# 3) New requirement added in E31
$ws.Range("E31").Value = "המערכת תיאפשר בחירת תווח זמנים להנפקת דוחות "

# 4) Apply yellow highlight fill to the requirement cells that got highlighted
$yellow = 65535

$bCells = @(2,3,14,15,16,17,22,23,24,25,26,27,28,29,30,31)
foreach ($r in $bCells) {
    $ws.Range("B$r").Interior.Color = $yellow
}

$eCells = @(2,3,4,5,6,8,9,10,11,12,13,14,15,16,17,18,19,25,27,29,31)
foreach ($r in $eCells) {
    $ws.Range("E$r").Interior.Color = $yellow
}

# 5) Sheet view changes
$ws.Activate()
$excel.ActiveWindow.Zoom = 141
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("E17").Select()
